# PB [23] - Funcionario Editar Stock
#
# This edit corrects two user-story rows that had accidentally been
# duplicated with a " NEW METHOD" suffix (rows 22/23), fills in the
# estimated/actual effort for the "alterar o stock dos produtos da loja"
# story (row 24, the "Funcionario Editar Stock" story), and moves three
# "Funcionario" backlog items from "To Do" into "Doing" with their
# Prioridade/Sprint filled in (rows 26-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (US 20) - drop the stray " NEW METHOD" suffix
$ws.Range("C22").Value = "quero apagar conta de um user"

# Row 23 (US 21) - drop the stray " NEW METHOD" suffix
$ws.Range("C23").Value = "quero editar os dados biométricos de um user"

# Row 24 (US 22) - "quero alterar o stock dos produtos da loja" now has an
# estimated and real effort of 20min
$ws.Range("F24").Value = "20min"
$ws.Range("G24").Value = "20min"

# Row 26 (US 24) - "quero gerir a lotação no ginásio" moves to Doing,
# priority 4, sprint 14
$ws.Range("E26").Value = 4
$ws.Range("I26").Value = "Doing"
$ws.Range("J26").Value = 14

# Row 27 (US 25) - "quero gerir os clientes" moves to Doing,
# priority 3, sprint 14
$ws.Range("E27").Value = 3
$ws.Range("I27").Value = "Doing"
$ws.Range("J27").Value = 14

# Row 28 (US 26) - "quero gerir as avaliações do ginásio" moves to Doing,
# priority 3, sprint 14
$ws.Range("E28").Value = 3
$ws.Range("I28").Value = "Doing"
$ws.Range("J28").Value = 14

# Row 29 (US 27) - "quero gerir a loja online" stays To Do, but priority and
# sprint are now filled in
$ws.Range("E29").Value = 3
$ws.Range("J29").Value = 14

$null = $ws.Range("C23").Select()
